$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = 578
$ws.Range("BG2").Value = 221
$ws.Range("E3").Value = 6.5238095238095
$ws.Range("BQ3").Value = 137
$ws.Range("AB4").Value = 352
$ws.Range("AF4").Value = 73.48066298342501
$ws.Range("AG4").Value = 181
$ws.Range("BG4").Value = 126
$ws.Range("BT4").Value = 48
$ws.Range("DG4").Value = 138
$ws.Range("AA5").Value = 0.35935256
$ws.Range("AB5").Value = 176
$ws.Range("AE5").Value = 60
$ws.Range("AF5").Value = 74.07407407407401
$ws.Range("AI5").Value = 41
$ws.Range("AJ5").Value = 29
$ws.Range("BG5").Value = 65
$ws.Range("BT5").Value = 21
$ws.Range("AA6").Value = 0.3150909
$ws.Range("AB7").Value = 380
$ws.Range("AA8").Value = 1.09980935
$ws.Range("AE8").Value = 286
$ws.Range("AF8").Value = 65.148063781321
$ws.Range("AH8").Value = 151
$ws.Range("AL8").Value = 29.310344827586
$ws.Range("AR8").Value = 70
$ws.Range("BA8").Value = 90
$ws.Range("BB8").Value = 50.847457627119
$ws.Range("BE8").Value = 25
$ws.Range("BF8").Value = 45.454545454545
$ws.Range("BG8").Value = 199
$ws.Range("BT8").Value = 153
$ws.Range("CZ8").Value = 58
$ws.Range("E12").Value = 6.8421052631579
$ws.Range("AV12").Value = 0
$ws.Range("BQ12").Value = 130
$ws.Range("CC12").Value = 101
$ws.Range("AB13").Value = 801
$ws.Range("BG13").Value = 182
$ws.Range("AB14").Value = 174
$ws.Range("E17").Value = 6.655
$ws.Range("AA17").Value = 2.32567956
$ws.Range("AE17").Value = 281
$ws.Range("AF17").Value = 77.410468319559
$ws.Range("AI17").Value = 159
$ws.Range("AJ17").Value = 88
$ws.Range("BG17").Value = 210
$ws.Range("BQ17").Value = 133.1
$ws.Range("BT17").Value = 82
$ws.Range("AS18").Value = 10
$ws.Range("BB18").Value = 63.114754098361
$ws.Range("BD18").Value = 60.550458715596
$ws.Range("CU18").Value = 90
$ws.Range("K20").Value = 3.038
$ws.Range("BA20").Value = 147
$ws.Range("BB20").Value = 63.636363636364
$ws.Range("BE20").Value = 101
$ws.Range("BF20").Value = 66.447368421053
$ws.Range("AA21").Value = 0.23373121
$ws.Range("AE21").Value = 92
$ws.Range("AF21").Value = 63.888888888889
$ws.Range("AG21").Value = 144
$ws.Range("AH21").Value = 49
$ws.Range("AO21").Value = 7
$ws.Range("AS21").Value = 0
$ws.Range("AT21").Value = 20
$ws.Range("BA21").Value = 23
$ws.Range("BB21").Value = 71.875
$ws.Range("BD21").Value = 66.666666666667
$ws.Range("BE21").Value = 11
$ws.Range("BF21").Value = 78.571428571429
$ws.Range("CU21").Value = 9
$ws.Range("DF21").Value = 60
$ws.Range("AB22").Value = 880
$ws.Range("CO22").Value = 26
$ws.Range("DL22").Value = -3.6812
